# Auto-generated Excel COM-interop script applying the Jenova_Profits.xlsx diff
# Updates market-price / profit columns (H-N) on specific rows across 7 of the 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8160.6
$ws.Range("J40").Value = 9920.200000000001
$ws.Range("L40").Value = 9920.200000000001
$ws.Range("N40").Value = -10270.2

$ws.Range("H64").Value = 10800.4
$ws.Range("J64").Value = 14333.333
$ws.Range("L64").Value = 14333.333
$ws.Range("N64").Value = -14829.333

$ws.Range("H67").Value = 10800.4
$ws.Range("J67").Value = 14333.333
$ws.Range("L67").Value = 14333.333
$ws.Range("N67").Value = -16049.333

$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -64992

$ws.Range("H98").Value = 1591.0714
$ws.Range("I98").Value = 1295.625
$ws.Range("J98").Value = 7500
$ws.Range("K98").Value = 1295.625
$ws.Range("L98").Value = 7500
$ws.Range("M98").Value = 202.375
$ws.Range("N98").Value = -10496

$ws.Range("H111").Value = 64431.883
$ws.Range("I111").Value = 95814
$ws.Range("K111").Value = 287442
$ws.Range("M111").Value = -284375

$ws.Range("H122").Value = 1591.0714
$ws.Range("I122").Value = 1295.625
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 3886.875
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -1436.875
$ws.Range("N122").Value = -27400

$ws.Range("H135").Value = 1054232.4
$ws.Range("I135").Value = 1334680.4
$ws.Range("K135").Value = 12012123.6
$ws.Range("M135").Value = -12009588.6

$ws.Range("H138").Value = 4278.275
$ws.Range("I138").Value = 1920.4
$ws.Range("J138").Value = 4822.4
$ws.Range("K138").Value = 5761.200000000001
$ws.Range("L138").Value = 14467.2
$ws.Range("M138").Value = -621.2000000000007
$ws.Range("N138").Value = -24747.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3796.6516
$ws.Range("I32").Value = 3660.276
$ws.Range("J32").Value = 9729
$ws.Range("K32").Value = 3660.276
$ws.Range("L32").Value = 9729
$ws.Range("M32").Value = -3373.276
$ws.Range("N32").Value = -10303

$ws.Range("H41").Value = 9000
$ws.Range("I41").Value = 9000
$ws.Range("K41").Value = 9000
$ws.Range("M41").Value = -8586

$ws.Range("H45").Value = 1931.2222
$ws.Range("J45").Value = 1999.3334
$ws.Range("L45").Value = 1999.3334
$ws.Range("N45").Value = -2753.3334

$ws.Range("H122").Value = 4291.905
$ws.Range("J122").Value = 6120.8335
$ws.Range("L122").Value = 18362.5005
$ws.Range("N122").Value = -23262.5005

$ws.Range("H132").Value = 5420.364
$ws.Range("I132").Value = 4962.4
$ws.Range("K132").Value = 14887.2
$ws.Range("M132").Value = -12357.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 169622.33
$ws.Range("J42").Value = 169622.33
$ws.Range("L42").Value = 169622.33
$ws.Range("N42").Value = -170278.33

$ws.Range("H43").Value = 223842
$ws.Range("J43").Value = 223842
$ws.Range("L43").Value = 223842
$ws.Range("N43").Value = -224204

$ws.Range("H134").Value = 63111.47
$ws.Range("I134").Value = 4084.2307
$ws.Range("K134").Value = 12252.6921
$ws.Range("M134").Value = -9717.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 494.6
$ws.Range("I7").Value = 511.42105
$ws.Range("J7").Value = 441.33334
$ws.Range("K7").Value = 511.42105
$ws.Range("L7").Value = 441.33334
$ws.Range("M7").Value = -398.42105
$ws.Range("N7").Value = -667.33334

$ws.Range("H22").Value = 249
$ws.Range("I22").Value = 110
$ws.Range("J22").Value = 341.66666
$ws.Range("K22").Value = 110
$ws.Range("L22").Value = 341.66666
$ws.Range("M22").Value = 240
$ws.Range("N22").Value = -1041.66666

$ws.Range("H31").Value = 42886.56
$ws.Range("J31").Value = 171834.83
$ws.Range("L31").Value = 171834.83
$ws.Range("N31").Value = -172424.83

$ws.Range("H34").Value = 42886.56
$ws.Range("J34").Value = 171834.83
$ws.Range("L34").Value = 171834.83
$ws.Range("N34").Value = -172238.83

$ws.Range("H86").Value = 7691.6
$ws.Range("I86").Value = 6701.5
$ws.Range("J86").Value = 8351.666999999999
$ws.Range("K86").Value = 6701.5
$ws.Range("L86").Value = 8351.666999999999
$ws.Range("M86").Value = -5578.5
$ws.Range("N86").Value = -10597.667

$ws.Range("H89").Value = 7691.6
$ws.Range("I89").Value = 6701.5
$ws.Range("J89").Value = 8351.666999999999
$ws.Range("K89").Value = 33507.5
$ws.Range("L89").Value = 41758.335
$ws.Range("M89").Value = -27891.5
$ws.Range("N89").Value = -52990.335

$ws.Range("H94").Value = 2975
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2975
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 2975
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -3877

$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 6333.3335
$ws.Range("I48").Value = 4000
$ws.Range("J48").Value = 7500
$ws.Range("K48").Value = 4000
$ws.Range("L48").Value = 7500
$ws.Range("M48").Value = -3515
$ws.Range("N48").Value = -8470

$ws.Range("H80").Value = 2865928.2
$ws.Range("I80").Value = 2008398.2
$ws.Range("K80").Value = 2008398.2
$ws.Range("M80").Value = -2007400.2

$ws.Range("H83").Value = 2865928.2
$ws.Range("I83").Value = 2008398.2
$ws.Range("K83").Value = 10041991
$ws.Range("M83").Value = -10036999

$ws.Range("H102").Value = 2089.3333
$ws.Range("I102").Value = 895.94116
$ws.Range("K102").Value = 895.94116
$ws.Range("M102").Value = 726.05884

$ws.Range("H132").Value = 274750
$ws.Range("I132").Value = 33000
$ws.Range("J132").Value = 1000000
$ws.Range("K132").Value = 99000
$ws.Range("L132").Value = 3000000
$ws.Range("M132").Value = -96470
$ws.Range("N132").Value = -3005060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1119691.2
$ws.Range("I7").Value = 17252
$ws.Range("J7").Value = 1434673.9
$ws.Range("K7").Value = 17252
$ws.Range("L7").Value = 1434673.9
$ws.Range("M7").Value = -17140
$ws.Range("N7").Value = -1434897.9

$ws.Range("H16").Value = 166667760
$ws.Range("I16").Value = 166667760
$ws.Range("K16").Value = 166667760
$ws.Range("M16").Value = -166667590

$ws.Range("H40").Value = 6254573.5
$ws.Range("I40").Value = 10003118
$ws.Range("K40").Value = 10003118
$ws.Range("M40").Value = -10002982

$ws.Range("H61").Value = 4306.04
$ws.Range("I61").Value = 3274.2942
$ws.Range("J61").Value = 6498.5
$ws.Range("K61").Value = 3274.2942
$ws.Range("L61").Value = 6498.5
$ws.Range("M61").Value = -3072.2942
$ws.Range("N61").Value = -6902.5

$ws.Range("H113").Value = 4306.04
$ws.Range("I113").Value = 3274.2942
$ws.Range("J113").Value = 6498.5
$ws.Range("K113").Value = 3274.2942
$ws.Range("L113").Value = 6498.5
$ws.Range("M113").Value = -1104.2942
$ws.Range("N113").Value = -10838.5

$ws.Range("H126").Value = 1119691.2
$ws.Range("I126").Value = 17252
$ws.Range("J126").Value = 1434673.9
$ws.Range("K126").Value = 51756
$ws.Range("L126").Value = 4304021.699999999
$ws.Range("M126").Value = -49286
$ws.Range("N126").Value = -4308961.699999999

$ws.Range("H132").Value = 2999.5
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 1058470
$ws.Range("I136").Value = 2227586
$ws.Range("K136").Value = 6682758
$ws.Range("M136").Value = -6680208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1996.8125
$ws.Range("I81").Value = 1805.625
$ws.Range("K81").Value = 3611.25
$ws.Range("M81").Value = -2550.25

$ws.Range("H84").Value = 1996.8125
$ws.Range("I84").Value = 1805.625
$ws.Range("K84").Value = 18056.25
$ws.Range("M84").Value = -12752.25

$ws.Range("H135").Value = 98549.60000000001
$ws.Range("J135").Value = 98549.60000000001
$ws.Range("L135").Value = 98549.60000000001
$ws.Range("N135").Value = -108689.6

$ws.Range("H136").Value = 8224343.5
$ws.Range("I136").Value = 9037064
$ws.Range("K136").Value = 27111192
$ws.Range("M136").Value = -27108642

$ws.Range("H138").Value = 333389630
$ws.Range("J138").Value = 333389630
$ws.Range("L138").Value = 333389630
$ws.Range("N138").Value = -333399910
